$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 2).Value = 0.5696291547463442
$ws.Cells.Item(4, 3).Value = 0.5640000000000001
$ws.Cells.Item(4, 4).Value = 0.6289259988206177
$ws.Cells.Item(4, 5).Value = 0.6075
$ws.Cells.Item(4, 6).Value = 0.6641157913352068
$ws.Cells.Item(4, 7).Value = 0.985
$ws.Cells.Item(4, 8).Value = 0.500948808950948
$ws.Cells.Item(4, 9).Value = 0.502
$ws.Cells.Item(4, 10).Value = 0.6091506016182076
$ws.Cells.Item(4, 11).Value = 0.65
$ws.Cells.Item(4, 12).Value = 0.5865631528644721
$ws.Cells.Item(4, 13).Value = 0.6110000000000001

$ws.Cells.Item(5, 2).Value = 0.6889796196133278
$ws.Cells.Item(5, 3).Value = 0.7180000000000001
$ws.Cells.Item(5, 4).Value = 0.674428656652411
$ws.Cells.Item(5, 5).Value = 0.6855
$ws.Cells.Item(5, 6).Value = 0.6473491909008886
$ws.Cells.Item(5, 7).Value = 0.9480000000000001
$ws.Cells.Item(5, 8).Value = 0.4918555057391506
$ws.Cells.Item(5, 9).Value = 0.487
$ws.Cells.Item(5, 10).Value = 0.6017323807093553
$ws.Cells.Item(5, 11).Value = 0.611
$ws.Cells.Item(5, 12).Value = 0.6032489817106115
$ws.Cells.Item(5, 13).Value = 0.6214999999999999

$ws.Cells.Item(6, 2).Value = 0.630258001061337
$ws.Cells.Item(6, 3).Value = 0.723
$ws.Cells.Item(6, 4).Value = 0.5801521608554855
$ws.Cells.Item(6, 5).Value = 0.593
$ws.Cells.Item(6, 6).Value = 0.663189793483707
$ws.Cells.Item(6, 7).Value = 0.9860000000000001
$ws.Cells.Item(6, 8).Value = 0.4996475575957275
$ws.Cells.Item(6, 9).Value = 0.4995
$ws.Cells.Item(6, 10).Value = 0.6302343615535737
$ws.Cells.Item(6, 11).Value = 0.6380000000000001
$ws.Cells.Item(6, 12).Value = 0.6362502519737273
$ws.Cells.Item(6, 13).Value = 0.6475

$ws.Cells.Item(7, 2).Value = 0.2154210467690116
$ws.Cells.Item(7, 3).Value = 0.176
$ws.Cells.Item(7, 4).Value = 0.3613725041177805
$ws.Cells.Item(7, 5).Value = 0.421
$ws.Cells.Item(7, 6).Value = 0.07849586080043762
$ws.Cells.Item(7, 7).Value = 0.102
$ws.Cells.Item(7, 8).Value = 0.2328125
$ws.Cells.Item(7, 9).Value = 0.4865
$ws.Cells.Item(7, 10).Value = 0.4032295676642556
$ws.Cells.Item(7, 11).Value = 0.4069999999999999
$ws.Cells.Item(7, 12).Value = 0.4033372780918357
$ws.Cells.Item(7, 13).Value = 0.4135

Write-Output "Updated B4:M7 with new classification results"